# Updated cryptos list on Wed Oct 11 07:50:00 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) figures for the cryptos table,
# and fix the swapped HuobiToken / TrustWalletToken rows (36 & 37).
#
# Price values in column D are plain text (e.g. "27.108.04", "6.60") even
# though some of them look numeric - Excel's automatic type inference would
# otherwise coerce "6.60" -> 6.6, dropping the trailing zero. Force those
# cells to Text format before assigning, then restore "Normal" style so no
# stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.108.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.560.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.20%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("E6").Value = '  -3.30%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.782.20'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.555.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("E14").Value = '  -2.91%  '
$ws.Range("E15").Value = '  -3.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.090.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("E19").Value = '  -2.03%  '
$ws.Range("E20").Value = '  -2.09%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("E23").Value = '  -4.78%  '
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("E31").Value = '  -2.52%  '
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.383.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.941'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.28%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("E38").Value = '  -1.98%  '
$ws.Range("E39").Value = '  -2.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.515'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.14%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.989'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.13%  '
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.695.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.23'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("E49").Value = '  -2.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0492'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("E51").Value = '  -0.04%  '
